$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text cells (B/C) to stay as text instead of being auto-converted
# to a date serial / number by Excel's input parser, then restore the
# default (unstyled) formatting so no stray style gets introduced.
$textRange = $ws.Range("B185:C186")
$textRange.NumberFormat = "@"

# Row 185: 2020-02-27
$ws.Cells.Item(185, 1).Value = 1582761600
$ws.Cells.Item(185, 2).Value = "2020-02-27"
$ws.Cells.Item(185, 3).Value = "6633"
$ws.Cells.Item(185, 4).Value = "LHI"
$ws.Cells.Item(185, 5).Value = 0.715
$ws.Cells.Item(185, 6).Value = 0.725
$ws.Cells.Item(185, 7).Value = 0.715
$ws.Cells.Item(185, 8).Value = 0.715
$ws.Cells.Item(185, 9).Value = 816500

# Row 186: 2020-02-28
$ws.Cells.Item(186, 1).Value = 1582848000
$ws.Cells.Item(186, 2).Value = "2020-02-28"
$ws.Cells.Item(186, 3).Value = "6633"
$ws.Cells.Item(186, 4).Value = "LHI"
$ws.Cells.Item(186, 5).Value = 0.71
$ws.Cells.Item(186, 6).Value = 0.71
$ws.Cells.Item(186, 7).Value = 0.6899999999999999
$ws.Cells.Item(186, 8).Value = 0.7
$ws.Cells.Item(186, 9).Value = 5056700

$textRange.ClearFormats()
